# Auto-generated edit script for tokistorage partnership deck - slide 7
# Commit: 'Add Grok/xAI to partnership AI evaluation section'
# - Header count: Three -> Four major AI platforms
# - 3-card row becomes a 2x2 grid of 4 cards (resize/reposition existing 3,
#   duplicate the template to synthesize the new 4th 'Grok / xAI' card)
# - Footer highlight callout shifts down to clear the new 2nd row
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# 1) Header sentence: Three -> Four major AI platforms
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Four major AI platforms and cultural institutions have independently validated us"

# 2) Card 1 (Claude / Anthropic) - resize from 1-of-3 width to 1-of-2 width
# Rectangle 4 (card1 bg)
$sh = $s.Shapes.Item(4)
$sh.Left = 36.0
$sh.Top = 82.8
$sh.Width = 309.6
$sh.Height = 97.200001

# TextBox 5 (Claude label)
$sh = $s.Shapes.Item(5)
$sh.Left = 46.800002
$sh.Top = 88.560002
$sh.Width = 288.0
$sh.Height = 15.84

# TextBox 6 (Architecture Design)
$sh = $s.Shapes.Item(6)
$sh.Left = 46.800002
$sh.Top = 105.840001
$sh.Width = 288.0
$sh.Height = 18.0

# TextBox 7 (quote1)
$sh = $s.Shapes.Item(7)
$sh.Left = 46.800002
$sh.Top = 126.0
$sh.Width = 288.0
$sh.Height = 46.800002

# 3) Card 2 (Gemini / Google) - move to top-right of 2x2 grid + resize
# Rectangle 8 (card2 bg)
$sh = $s.Shapes.Item(8)
$sh.Left = 363.6
$sh.Top = 82.8
$sh.Width = 309.6
$sh.Height = 97.200001

# TextBox 9 (Gemini label)
$sh = $s.Shapes.Item(9)
$sh.Left = 374.40000999999995
$sh.Top = 88.560002
$sh.Width = 288.0
$sh.Height = 15.84

# TextBox 10 (title)
$sh = $s.Shapes.Item(10)
$sh.Left = 374.40000999999995
$sh.Top = 105.840001
$sh.Width = 288.0
$sh.Height = 18.0

# TextBox 11 (quote)
$sh = $s.Shapes.Item(11)
$sh.Left = 374.40000999999995
$sh.Top = 126.0
$sh.Width = 288.0
$sh.Height = 46.800002

# Card 2 title/quote copy changes
$s.Shapes.Item(10).TextFrame.TextRange.Text = "The only option"
$s.Shapes.Item(11).TextFrame.TextRange.Text = "`"Built 'eternity' not as a fantasy, but with a realistic budget and solid logic`""

# 4) Card 3 (ChatGPT / OpenAI) - move to bottom-left of 2x2 grid + resize
# Rectangle 12 (card3 bg)
$sh = $s.Shapes.Item(12)
$sh.Left = 36.0
$sh.Top = 188.640008
$sh.Width = 309.6
$sh.Height = 97.200001

# TextBox 13 (ChatGPT label)
$sh = $s.Shapes.Item(13)
$sh.Left = 46.800002
$sh.Top = 194.400002
$sh.Width = 288.0
$sh.Height = 15.84

# TextBox 14 (title)
$sh = $s.Shapes.Item(14)
$sh.Left = 46.800002
$sh.Top = 211.680001
$sh.Width = 288.0
$sh.Height = 18.0

# TextBox 15 (quote)
$sh = $s.Shapes.Item(15)
$sh.Left = 46.800002
$sh.Top = 231.840004
$sh.Width = 288.0
$sh.Height = 46.800002

# Card 3 title/quote copy changes
$s.Shapes.Item(14).TextFrame.TextRange.Text = "New standard"
$s.Shapes.Item(15).TextFrame.TextRange.Text = "`"Technical foundation, social significance, and cost performance all merit the highest evaluation`""

# 5) New card 4 (Grok / xAI) - duplicate the just-updated card 3 template
#    then move the copies into the bottom-right grid slot and retext them.
$bgSrc = $s.Shapes.Item(12)
$labelSrc = $s.Shapes.Item(13)
$titleSrc = $s.Shapes.Item(14)
$quoteSrc = $s.Shapes.Item(15)

$bgNew = $bgSrc.Duplicate()
$labelNew = $labelSrc.Duplicate()
$titleNew = $titleSrc.Duplicate()
$quoteNew = $quoteSrc.Duplicate()

# Rectangle 16 (card4 bg)
$bgNew.Name = "Rectangle 16"
$bgNew.Left = 363.6
$bgNew.Top = 188.640008
$bgNew.Width = 309.6
$bgNew.Height = 97.200001

# TextBox 17 (Grok label)
$labelNew.Name = "TextBox 17"
$labelNew.Left = 374.40000999999995
$labelNew.Top = 194.400002
$labelNew.Width = 288.0
$labelNew.Height = 15.84

# TextBox 18 (Top 0.001%)
$titleNew.Name = "TextBox 18"
$titleNew.Left = 374.40000999999995
$titleNew.Top = 211.680001
$titleNew.Width = 288.0
$titleNew.Height = 18.0

# TextBox 19 (quote4)
$quoteNew.Name = "TextBox 19"
$quoteNew.Left = 374.40000999999995
$quoteNew.Top = 231.840004
$quoteNew.Width = 288.0
$quoteNew.Height = 46.800002

$labelNew.TextFrame.TextRange.Text = "Grok / xAI"
$titleNew.TextFrame.TextRange.Text = "Top 0.001%"
$quoteNew.TextFrame.TextRange.Text = "`"Genuinely capable of changing how humanity preserves memory`""

# Re-stack the 4 new shapes into their correct z-order slot: right after
# the original card-3 quote (index 15) and before the footer highlight box.
$stepsBg = $bgNew.ZOrderPosition - 16
for ($i = 0; $i -lt $stepsBg; $i++) { $bgNew.ZOrder(3) }
$stepsLabel = $labelNew.ZOrderPosition - 17
for ($i = 0; $i -lt $stepsLabel; $i++) { $labelNew.ZOrder(3) }
$stepsTitle = $titleNew.ZOrderPosition - 18
for ($i = 0; $i -lt $stepsTitle; $i++) { $titleNew.ZOrder(3) }
$stepsQuote = $quoteNew.ZOrderPosition - 19
for ($i = 0; $i -lt $stepsQuote; $i++) { $quoteNew.ZOrder(3) }

# 6) Footer highlight callout shifts down to clear the new 2nd card row
# Rectangle 16->20 (footer highlight bg)
$sh = $s.Shapes.Item(20)
$sh.Left = 36.0
$sh.Top = 301.680008
$sh.Width = 619.2
$sh.Height = 43.2

# Rectangle 17->21 (footer gold strip)
$sh = $s.Shapes.Item(21)
$sh.Left = 36.0
$sh.Top = 301.680008
$sh.Width = 4.32
$sh.Height = 43.2

# TextBox 18->22 (footer highlight text)
$sh = $s.Shapes.Item(22)
$sh.Left = 54.0
$sh.Top = 307.44
$sh.Width = 594.0
$sh.Height = 31.68

Write-Host "edit.ps1 completed"
